# Auto-generated from OOXML diff: BRVM "Recommandations" + "Top_YTD" automatic update
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# Full replacement data for Recommandations!A2:G54 (row order/content changed by the BRVM refresh)
$recoData = @(
    @("CFAO MOTORS CI", 0, 5, 4700, 940, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 5, 4650, 935, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 10, 4200.86, 109.01, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 5, 3470, 695, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 5, 3310, 655, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 5, 3117.28, 648.9400000000001, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 5, 2004.06, 405.43, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 5, 1801.53, 364.71, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 5, 1692.85, 337.46, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI", 0, 1, 990, 990, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 5, 710.1900000000001, 142.55, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 5, 682.09, 137.5, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 5, 665.54, 134.61, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 5, 654.0700000000001, 132.29, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 5, 631.23, 129.6, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 5, 533.26, 107.45, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 5, 476.11, 95.87, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE              (**)", 0, 2, 452.63, 224.69, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE             (**)", 0, 2, 382.52, 188.06, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                    (**)", 0, 1, 218.25, 218.25, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                 (**)", 0, 1, 216.31, 216.31, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 192.72, 192.72, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                 (**)", 0, 1, 192.29, 192.29, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 1, 192.05, 192.05, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                       (**)", 0, 1, 191.44, 191.44, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                (**)", 0, 1, 190.98, 190.98, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE            (**)", 0, 1, 186.07, 186.07, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 3, 0, 21.1, 7.07, "🟢 Achat", "✅ Renforcer"),
    @("SAFCA CI (SAFC)", 3, 0, 19.04, 6.49, "🟢 Achat", "✅ Renforcer"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 3, 0, 15.31, 7.46, "🟢 Achat", "✅ Renforcer"),
    @("UNIWAX CI (UNXC)", 2, 0, 14.15, 7.08, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 2, 1, 11, -2.65, "🟡 Observer", "👀 À surveiller"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 1, 0, 7.49, 7.49, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA CI (BOAC)", 1, 0, 7.44, 7.44, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 2, 1, 6.2, 5.88, "🟡 Observer", "👀 À surveiller"),
    @("SUCRIVOIRE (SCRC)", 2, 2, 4.73, -1.82, "🟡 Observer", "👀 À surveiller"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 1, 0, 3.9, 3.9, "🟡 Observer", "➖ Neutre"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 1, 0, 1.72, 1.72, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 1, 1, 1.38, 4.87, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 5, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 1, 1, -0.27, 3.7, "🟡 Observer", "👀 À surveiller"),
    @("FILTISAC CI (FTSC)", 0, 1, -0.63, -0.63, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 0, 1, -0.65, -0.65, "🟡 Observer", "➖ Neutre"),
    @("PALM CI (PALC)", 0, 1, -1.1, -1.1, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 0, 1, -1.15, -1.15, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI (SIVC)", 0, 1, -1.44, -1.44, "🟡 Observer", "➖ Neutre"),
    @("SMB CI (SMBC)", 1, 1, -2.94, 2.17, "🟡 Observer", "👀 À surveiller"),
    @("SOGB CI (SOGC)", 0, 1, -3.06, -3.06, "🟡 Observer", "➖ Neutre"),
    @("CORIS BANK INTERNATIONAL (CBIBF)", 1, 2, -3.94, 4.52, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA ML (BOAM)", 0, 1, -4.88, -4.88, "🟡 Observer", "➖ Neutre"),
    @("SICOR CI (SICC)", 0, 1, -7.14, -7.14, "🟡 Observer", "➖ Neutre"),
    @("ORAGROUP TOGO (ORGT)", 0, 2, -12.45, -6.17, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 0, 5, -37.33, -7.48, "🔴 Vente", "⚠️ Risque de décrochage")
)

for ($i = 0; $i -lt $recoData.Count; $i++) {
    $r = 2 + $i
    $vals = $recoData[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $wsReco.Cells.Item($r, 1 + $j).Value = $vals[$j]
    }
}

# dimension now extends to row 54 (was 48) - handled automatically by the extra rows above

# Top_YTD!B2:B11 refreshed totals
$ytdData = @(
    158811865.19,
    12163032.8,
    11563967.94,
    3155474.42,
    2568786.1,
    1978657.87,
    314916.51,
    206518.88,
    162149.21,
    8206.43
)
for ($i = 0; $i -lt $ytdData.Count; $i++) {
    $wsYtd.Cells.Item(2 + $i, 2).Value = $ytdData[$i]
}
